# Slide 4 ("LA METODOLOGÍA DE INVESTIGACIÓN") - last paragraph of the single
# textbox ("CuadroTexto 2") originally read:
#   "...la clínica PLAyANA ha decidió crear una plataforma web, ..."
# and must become:
#   "...la clínica PLAyANA creará una plataforma web, ..."
# split, per the target OOXML, into three runs:
#   1) "...la clínica PLAyANA crear"
#   2) "á"
#   3) " una plataforma web, ..."

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(4)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(5)

$newPrefix = "El país de Guinea Ecuatorial carece de clínicas especializadas en servicios oftalmólogos, y teniendo en cuenta que en Guinea Ecuatorial el acceso a internet y el uso de teléfonos móviles es cada vez mayor, la clínica PLAyANA crear"
$oldPrefixLen = $newPrefix.Length + "ha decidió ".Length

# 1) Rewrite the paragraph's leading characters - everything up to and
#    including "...PLAyANA ha decidió crear" - dropping "ha decidió ", so
#    the run keeps going right up to "crear" with no extra split introduced.
$prefix = $para.Characters(1, $oldPrefixLen)
$prefix.Text = $newPrefix

# 2) The single space that used to sit right after "crear" becomes its own
#    run containing just "á".
$accentPos = $newPrefix.Length + 1
$accentRun = $para.Characters($accentPos, 1)
$accentRun.Text = "á"

# 3) Restore the space in front of "una plataforma..." as part of that
#    trailing run (merges with it, producing a third, separate run).
$tailPos = $accentPos + 1
$tail = $para.Characters($tailPos, $para.Length - $tailPos + 1)
[void]$tail.InsertBefore(" ")
